# Weekly fruit/vegetable price update: a new observation for
# "Comercializadora del Agro de Limarí - Arveja Verde" is inserted as the
# new row 67 (date 2023-07-13 / serial 45120), pushing the previously
# existing rows 67-93 down to 68-94.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 67, shifting rows 67:93 -> 68:94
# (and bringing the A1:R93 dimension to A1:R94).
$ws.Rows.Item(67).Insert()

# Populate the newly inserted row with this week's record.
$ws.Cells.Item(67, 1).Value = 2
$ws.Cells.Item(67, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(67, 3).Value = "Coquimbo"
$ws.Cells.Item(67, 4).Value = 45120
$ws.Cells.Item(67, 5).Value = 4
$ws.Cells.Item(67, 6).Value = 100112022
$ws.Cells.Item(67, 7).Value = "Arveja Verde"
$ws.Cells.Item(67, 8).Value = "Perfection"
$ws.Cells.Item(67, 9).Value = "Primera"
$ws.Cells.Item(67, 10).Value = 1100
$ws.Cells.Item(67, 11).Value = 25000
$ws.Cells.Item(67, 12).Value = 27000
$ws.Cells.Item(67, 13).Value = 26000
$ws.Cells.Item(67, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(67, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(67, 16).Value = 1040
$ws.Cells.Item(67, 17).Value = 25
$ws.Cells.Item(67, 18).Value = "Hortaliza"
